$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FMICIC")
$v = $ws.Range("A17").Value
Write-Output ($v.GetType())
Write-Output ([string]$v)
